# Update filtered_output.xlsx: refresh the "Neg_Change" and "Pos_Change"
# sheets with the latest filtered market data (new symbols/rows from the
# upstream scrape), per the GitHub Actions scheduled refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Neg_Change  (dimension A1:I11 -> A1:I13)
# ---------------------------------------------------------------------
$wsNeg = $wb.Worksheets.Item("Neg_Change")

$negData = @(
  @("KOTAKBANK",  2090.6,  2107,    2080,    2095.6,  2060068,  5030328,   -0.5904704424840687),
  @("ADANIENT",   2462,    2478.6,  2434.3,  2439,    1245215,  2672697,   -0.5340979542387334),
  @("ADANIGREEN", 1092,    1097.2,  1075,    1077.2,  1414894,  2888868,   -0.5102254585533157),
  @("IRFC",       122.5,   122.74,  120.7,   120.92,  4555639,  9945324,   -0.5419315650249303),
  @("EXIDEIND",   382.2,   383.45,  378.95,  381.9,   1771009,  3658302,   -0.5158931657364537),
  @("BANKINDIA",  147.41,  148.7,   146.21,  146.99,  5679286,  11822101,  -0.5196043410557903),
  @("POLICYBZR",  1815,    1831,    1789,    1803,    1391840,  3081344,   -0.5483010011215885),
  @("INDUSINDBK", 854.05,  858.5,   845.05,  847.8,   1958051,  4019335,   -0.5128420497420593),
  @("YESBANK",    23.19,   23.25,   22.95,   22.97,   74216797, 156830924, -0.5267719203133688),
  @("MUTHOOTFIN", 3680.1,  3748.1,  3680.1,  3692.1,  525413,   1176503,   -0.5534112535199656),
  @("CESC",       175.61,  176,     172.82,  175.25,  983855,   2442534,   -0.5971990563897984),
  @("ANGELONE",   2839,    2867,    2785,    2790.9,  744032,   1513166,   -0.508294529483216)
)

for ($i = 0; $i -lt $negData.Length; $i++) {
  $r = $i + 2
  $row = $negData[$i]
  $symbol = $row[0]
  $wsNeg.Cells.Item($r, 1).Value = $symbol
  $wsNeg.Cells.Item($r, 2).Value = $row[1]
  $wsNeg.Cells.Item($r, 3).Value = $row[2]
  $wsNeg.Cells.Item($r, 4).Value = $row[3]
  $wsNeg.Cells.Item($r, 5).Value = $row[4]
  $wsNeg.Cells.Item($r, 6).Value = $row[5]
  $wsNeg.Cells.Item($r, 7).Value = $row[6]
  $wsNeg.Cells.Item($r, 8).Value = $row[7]
  $wsNeg.Cells.Item($r, 9).Value = $symbol
}

# ---------------------------------------------------------------------
# Sheet 2: Pos_Change  (dimension A1:I7 -> A1:I16)
# ---------------------------------------------------------------------
$wsPos = $wb.Worksheets.Item("Pos_Change")

$posData = @(
  @("ASIANPAINT", 2886,    2916,    2862.1,  2902.2,  1819066,  1230683,  0.4780946840087983),
  @("LT",         4020,    4026.8,  3978.9,  3998.9,  1748928,  1206133,  0.4500291427230662),
  @("INDIGO",     5876,    5888.5,  5733.5,  5758.5,  906483,   608336,   0.4901025091396859),
  @("TECHM",      1444,    1451.1,  1419.5,  1420.7,  1553597,  999416,   0.5545048308211996),
  @("ADANIENSOL", 1023.8,  1036,    1017.75, 1024.25, 1175321,  837345,   0.4036281341621434),
  @("AMBUJACEM",  568,     568,     556.75,  558,     904843,   618154,   0.4637824878590125),
  @("LTIM",       5848,    5862.5,  5749,    5762.5,  160191,   108410,   0.4776404390738861),
  @("BHEL",       284.95,  290.2,   282.9,   289.2,   11332882, 7910619,  0.4326163350807313),
  @("BSE",        2811.6,  2872.5,  2787.9,  2826.3,  5402562,  3818270,  0.4149240362782098),
  @("PATANJALI",  593,     594.95,  582,     586.05,  2112668,  1391118,  0.5186835336757917),
  @("PHOENIXLTD", 1745.9,  1749,    1717.7,  1729.9,  287932,   191049,   0.5071107412234558),
  @("COFORGE",    1790,    1795,    1771.5,  1775,    785781,   518247,   0.5162287480680062),
  @("OBEROIRLTY", 1745,    1753,    1705.4,  1707,    233440,   150217,   0.5540185198745814),
  @("NCC",        185,     185.12,  181.4,   182.18,  3443281,  2410328,  0.4285528774507038),
  @("PPLPHARMA",  195.47,  195.88,  191.6,   191.95,  2356015,  1619560,  0.4547253575045074)
)

for ($i = 0; $i -lt $posData.Length; $i++) {
  $r = $i + 2
  $row = $posData[$i]
  $symbol = $row[0]
  $wsPos.Cells.Item($r, 1).Value = $symbol
  $wsPos.Cells.Item($r, 2).Value = $row[1]
  $wsPos.Cells.Item($r, 3).Value = $row[2]
  $wsPos.Cells.Item($r, 4).Value = $row[3]
  $wsPos.Cells.Item($r, 5).Value = $row[4]
  $wsPos.Cells.Item($r, 6).Value = $row[5]
  $wsPos.Cells.Item($r, 7).Value = $row[6]
  $wsPos.Cells.Item($r, 8).Value = $row[7]
  $wsPos.Cells.Item($r, 9).Value = $symbol
}
